$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.336.95"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").Value = "3.303.18"
$ws.Range("E3").Value = "  -4.60%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.57"
$ws.Range("E5").Value = "  -6.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "647.07"
$ws.Range("E6").Value = "  -3.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.32"
$ws.Range("E7").Value = "  -16.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.405"
$ws.Range("E8").Value = "  -11.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.955"
$ws.Range("E10").Value = "  -15.99%  "

$ws.Range("D11").Value = "3.296.67"
$ws.Range("E11").Value = "  -4.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.202"
$ws.Range("E12").Value = "  -7.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.12"
$ws.Range("E13").Value = "  -10.28%  "

$ws.Range("D14").Value = "96.402.69"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.89"
$ws.Range("E15").Value = "  -5.83%  "

$ws.Range("D16").Value = "3.910.36"
$ws.Range("E16").Value = "  -5.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000246"
$ws.Range("E17").Value = "  -9.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.36"
$ws.Range("E18").Value = "  -5.72%  "

$ws.Range("D19").Value = "3.301.68"
$ws.Range("E19").Value = "  -4.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.48"
$ws.Range("E20").Value = "  -6.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.21"
$ws.Range("E21").Value = "  -7.26%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.20"
$ws.Range("E22").Value = "  -7.13%  "

$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.28"
$ws.Range("E23").Value = "  -8.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.451"
$ws.Range("E24").Value = "  -11.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000194"
$ws.Range("E25").Value = "  -10.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.32"
$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "91.26"
$ws.Range("E27").Value = "  -12.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.78"
$ws.Range("E28").Value = "  -10.50%  "

$ws.Range("D29").Value = "3.486.39"
$ws.Range("E29").Value = "  -4.67%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.138"
$ws.Range("E31").Value = "  -11.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.54"
$ws.Range("E32").Value = "  -11.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.182"
$ws.Range("E33").Value = "  -8.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  +6.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.531"
$ws.Range("E36").Value = "  -11.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.46"
$ws.Range("E37").Value = "  -10.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.37"
$ws.Range("E40").Value = "  -9.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.147"
$ws.Range("E41").Value = "  -8.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "494.09"
$ws.Range("E42").Value = "  -8.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.46"
$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.65"
$ws.Range("E44").Value = "  -2.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.805"
$ws.Range("E45").Value = "  -7.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0397"
$ws.Range("E46").Value = "  -10.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.17"
$ws.Range("E47").Value = "  -5.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.27"
$ws.Range("E48").Value = "  -2.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.58"
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.18"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.09"
$ws.Range("E51").Value = "  -11.01%  "
